$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.4
$ws.Range("I2").Value = 2.35
$ws.Range("L2").Value = 3.1
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("W2").Value = 9
$ws.Range("AI2").Value = 10
$ws.Range("AJ2").Value = 9.5
$ws.Range("AK2").Value = 21
$ws.Range("AX2").Value = 13
$ws.Range("AZ2").Value = 41
$ws.Range("BA2").Value = 67

# Row 4 updates
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6
